# Update the "想去人数" (interested count) values on the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 41
$ws1.Range("F4").Value = 2106
$ws1.Range("F5").Value = 172
$ws1.Range("F6").Value = 358

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 41
$ws4.Range("F4").Value = 2106
$ws4.Range("F5").Value = 172
$ws4.Range("F7").Value = 358
